$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.580.47'
$ws.Range('E2').Value = '  +1.88%  '
$ws.Range('D3').Value = '2.987.49'
$ws.Range('E3').Value = '  +2.64%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('E5').Value = '  +2.27%  '
$ws.Range('D6').Value = "'104.15"
$ws.Range('E6').Value = '  +5.04%  '
$ws.Range('E7').Value = '  +2.60%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').Value = "'0.596"
$ws.Range('E9').Value = '  +3.11%  '
$ws.Range('D10').Value = "'36.82"
$ws.Range('E10').Value = '  +3.57%  '
$ws.Range('E11').Value = '  -0.33%  '
$ws.Range('E12').Value = '  +2.63%  '
$ws.Range('D13').Value = '3.469.10'
$ws.Range('E13').Value = '  +2.86%  '
$ws.Range('D14').Value = "'18.49"
$ws.Range('E14').Value = '  +3.82%  '
$ws.Range('D15').Value = "'7.82"
$ws.Range('E15').Value = '  +4.61%  '
$ws.Range('D16').Value = '2.983.32'
$ws.Range('E16').Value = '  +2.57%  '
$ws.Range('D17').Value = "'11.27"
$ws.Range('E17').Value = '  +2.38%  '
$ws.Range('D18').Value = "'0.996"
$ws.Range('E18').Value = '  +1.83%  '
$ws.Range('D19').Value = '51.597.94'
$ws.Range('E19').Value = '  +1.96%  '
$ws.Range('E20').Value = '  +1.85%  '
$ws.Range('E21').Value = '  +2.62%  '
$ws.Range('D22').Value = '0.0₃0966'
$ws.Range('E22').Value = '  +2.26%  '
$ws.Range('D23').Value = "'70.43"
$ws.Range('E23').Value = '  +2.71%  '
$ws.Range('D24').Value = "'267.63"
$ws.Range('E24').Value = '  +1.64%  '
$ws.Range('E25').Value = '  +4.02%  '
$ws.Range('D26').Value = "'7.99"
$ws.Range('E26').Value = '  +0.69%  '
$ws.Range('D27').Value = "'0.170"
$ws.Range('E27').Value = '  +5.39%  '
$ws.Range('D28').Value = "'7.22"
$ws.Range('E28').Value = '  -0.80%  '
$ws.Range('E29').Value = '  +0.04%  '
$ws.Range('D30').Value = "'26.15"
$ws.Range('E30').Value = '  +3.35%  '
$ws.Range('E31').Value = '  +1.70%  '
$ws.Range('E32').Value = '  +5.02%  '
$ws.Range('D33').Value = "'34.60"
$ws.Range('E33').Value = '  +5.84%  '
$ws.Range('D34').Value = "'51.42"
$ws.Range('E34').Value = '  +1.38%  '
$ws.Range('E35').Value = '  +0.67%  '
$ws.Range('D36').Value = "'0.0445"
$ws.Range('E36').Value = '  +2.79%  '
$ws.Range('E37').Value = '  -0.05%  '
$ws.Range('E38').Value = '  +8.63%  '
$ws.Range('E39').Value = '  +4.26%  '
$ws.Range('D40').Value = "'2.58"
$ws.Range('E40').Value = '  +6.66%  '
$ws.Range('E41').Value = '  +2.27%  '
$ws.Range('B42').Value = 'ARBITRUM'
$ws.Range('C42').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D42').Value = "'1.84"
$ws.Range('E42').Value = '  +3.84%  '
$ws.Range('B43').Value = 'NEARProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D43').Value = "'3.86"
$ws.Range('E43').Value = '  +16.04%  '
$ws.Range('D44').Value = "'127.28"
$ws.Range('E44').Value = '  +7.14%  '
$ws.Range('D45').Value = "'21.39"
$ws.Range('E45').Value = '  +3.36%  '
$ws.Range('E46').Value = '  -0.32%  '
$ws.Range('D47').Value = "'0.272"
$ws.Range('E47').Value = '  +2.84%  '
$ws.Range('E48').Value = '  +0.66%  '
$ws.Range('D49').Value = '2.034.14'
$ws.Range('E49').Value = '  +3.16%  '
$ws.Range('D50').Value = '3.283.51'
$ws.Range('E50').Value = '  +2.30%  '
$ws.Range('B51').Value = 'BEAM'
$ws.Range('C51').Value = 'https://coinranking.com/coin/cYYMfXF4u+beam-beam'
$ws.Range('D51').Value = "'0.0331"
$ws.Range('E51').Value = '  +2.96%  '

# Reset style so the forced-text cells keep default (no quote-prefix) styling
$ws.Range('D6').Style = 'Normal'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D51').Style = 'Normal'
